# Apply the "testing different values" edit described by the commit.
# Core semantic change: update assumed probabilities on the
# `potential_preg_untrt` sheet (column C, rows 9-17, skipping row 12).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("potential_preg_untrt")

$ws.Range("C9").Value  = 0.05
$ws.Range("C10").Value = 0.02
$ws.Range("C11").Value = 0.02
$ws.Range("C13").Value = 0.005
$ws.Range("C14").Value = 0.004
$ws.Range("C15").Value = 0.004
$ws.Range("C16").Value = 0.004
$ws.Range("C17").Value = 0.004

# Reflect the reviewer's final selection/active sheet in the saved view.
$ws.Range("C2:C21").Select()
$ws.Activate()

$simParams = $wb.Worksheets.Item("SimParameters")
$simParams.Range("B15").Select()
